$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 469.75
$ws.Range("J9").Value = 99.5
$ws.Range("L9").Value = 99.5
$ws.Range("N9").Value = -437.5

$ws.Range("H12").Value = 289.8
$ws.Range("J12").Value = 288.33334
$ws.Range("L12").Value = 288.33334
$ws.Range("N12").Value = -628.33334

$ws.Range("H18").Value = 998.1875
$ws.Range("I18").Value = 937.75
$ws.Range("K18").Value = 937.75
$ws.Range("M18").Value = -653.75

$ws.Range("H28").Value = 1829.125
$ws.Range("J28").Value = 446.25
$ws.Range("L28").Value = 446.25
$ws.Range("N28").Value = -1416.25

$ws.Range("H41").Value = 726.2
$ws.Range("I41").Value = 657.75
$ws.Range("J41").Value = 1000
$ws.Range("K41").Value = 657.75
$ws.Range("L41").Value = 1000
$ws.Range("M41").Value = -217.75
$ws.Range("N41").Value = -1880

$ws.Range("H43").Value = 1334.7
$ws.Range("I43").Value = 1341.8334
$ws.Range("K43").Value = 1341.8334
$ws.Range("M43").Value = -1272.8334

$ws.Range("H70").Value = 2678.4285
$ws.Range("J70").Value = 2437.25
$ws.Range("L70").Value = 7311.75
$ws.Range("N70").Value = -7851.75

$ws.Range("H73").Value = 2678.4285
$ws.Range("J73").Value = 2437.25
$ws.Range("L73").Value = 7311.75
$ws.Range("N73").Value = -9183.75

$ws.Range("H96").Value = 8626.286
$ws.Range("I96").Value = 12252.556
$ws.Range("J96").Value = 2099
$ws.Range("K96").Value = 36757.66800000001
$ws.Range("L96").Value = 6297
$ws.Range("M96").Value = -35384.66800000001
$ws.Range("N96").Value = -9043

$ws.Range("H125").Value = 12516
$ws.Range("I125").Value = 10032
$ws.Range("J125").Value = 15000
$ws.Range("K125").Value = 90288
$ws.Range("L125").Value = 135000
$ws.Range("M125").Value = -87828
$ws.Range("N125").Value = -139920

$ws.Range("H129").Value = 2552.5454
$ws.Range("I129").Value = 696.5
$ws.Range("K129").Value = 2089.5
$ws.Range("M129").Value = 2910.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6422188
$ws.Range("I32").Value = 6369660
$ws.Range("K32").Value = 6369660
$ws.Range("M32").Value = -6369373

$ws.Range("H53").Value = 11990
$ws.Range("J53").Value = 11990
$ws.Range("L53").Value = 11990
$ws.Range("N53").Value = -13354

$ws.Range("H74").Value = 2166.6667
$ws.Range("I74").Value = 1750
$ws.Range("K74").Value = 1750
$ws.Range("M74").Value = -876

$ws.Range("H77").Value = 2166.6667
$ws.Range("I77").Value = 1750
$ws.Range("K77").Value = 8750
$ws.Range("M77").Value = -4382

$ws.Range("H97").Value = 847.1429000000001
$ws.Range("I97").Value = 688.7273
$ws.Range("K97").Value = 688.7273
$ws.Range("M97").Value = -192.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 61000
$ws.Range("I35").Value = 61000
$ws.Range("K35").Value = 61000
$ws.Range("M35").Value = -60690

$ws.Range("H105").Value = 2699.75
$ws.Range("I105").Value = 2999.5
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 2999.5
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -1252.5
$ws.Range("N105").Value = -5894

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 26438
$ws.Range("J74").Value = 26438
$ws.Range("L74").Value = 26438
$ws.Range("N74").Value = -28186

$ws.Range("H77").Value = 26438
$ws.Range("J77").Value = 26438
$ws.Range("L77").Value = 79314
$ws.Range("N77").Value = -88050

$ws.Range("H105").Value = 3227.5
$ws.Range("I105").Value = 2473.8333
$ws.Range("J105").Value = 3981.1667
$ws.Range("K105").Value = 2473.8333
$ws.Range("L105").Value = 3981.1667
$ws.Range("M105").Value = -726.8332999999998
$ws.Range("N105").Value = -7475.1667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 315773.28
$ws.Range("I2").Value = 220053.2
$ws.Range("J2").Value = 555073.5
$ws.Range("K2").Value = 1320319.2
$ws.Range("L2").Value = 3330441
$ws.Range("M2").Value = -1320206.2
$ws.Range("N2").Value = -3330667

$ws.Range("H4").Value = 17215112
$ws.Range("I4").Value = 20017632
$ws.Range("K4").Value = 60052896
$ws.Range("M4").Value = -60052784

$ws.Range("H6").Value = 310.14285
$ws.Range("I6").Value = 320.16666
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 960.4999799999999
$ws.Range("L6").Value = 750
$ws.Range("M6").Value = -847.4999799999999
$ws.Range("N6").Value = -976

$ws.Range("H40").Value = 528.55554
$ws.Range("I40").Value = 159.4
$ws.Range("K40").Value = 637.6
$ws.Range("M40").Value = -568.6

$ws.Range("H81").Value = 2119.6667
$ws.Range("J81").Value = 2119.6667
$ws.Range("L81").Value = 6359.000100000001
$ws.Range("N81").Value = -8605.000100000001

$ws.Range("H84").Value = 2119.6667
$ws.Range("J84").Value = 2119.6667
$ws.Range("L84").Value = 19077.0003
$ws.Range("N84").Value = -30309.0003

$ws.Range("H95").Value = 8875
$ws.Range("J95").Value = 8875
$ws.Range("L95").Value = 26625
$ws.Range("N95").Value = -30743

$ws.Range("H106").Value = 16971
$ws.Range("I106").Value = 9399.5
$ws.Range("K106").Value = 28198.5
$ws.Range("M106").Value = -27252.5

$ws.Range("H126").Value = 4000
$ws.Range("J126").Value = 4000
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -21880

$ws.Range("H132").Value = 1733.3334
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1638875.1
$ws.Range("J11").Value = 514999.66
$ws.Range("L11").Value = 514999.66
$ws.Range("N11").Value = -515277.66

$ws.Range("H52").Value = 30030
$ws.Range("I52").Value = 30030
$ws.Range("K52").Value = 30030
$ws.Range("M52").Value = -29771

$ws.Range("H107").Value = 1669.4166
$ws.Range("I107").Value = 316.125
$ws.Range("J107").Value = 4376
$ws.Range("K107").Value = 316.125
$ws.Range("L107").Value = 4376
$ws.Range("M107").Value = 1603.875
$ws.Range("N107").Value = -8216

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6617.3335
$ws.Range("I7").Value = 2877.8333
$ws.Range("K7").Value = 2877.8333
$ws.Range("M7").Value = -2765.8333

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H23").Value = 40006
$ws.Range("I23").Value = 40006
$ws.Range("K23").Value = 40006
$ws.Range("M23").Value = -39776

$ws.Range("H40").Value = 4570.143
$ws.Range("I40").Value = 4298
$ws.Range("J40").Value = 4679
$ws.Range("K40").Value = 4298
$ws.Range("L40").Value = 4679
$ws.Range("M40").Value = -4162
$ws.Range("N40").Value = -4951

$ws.Range("H46").Value = 3185.5
$ws.Range("J46").Value = 3585
$ws.Range("L46").Value = 3585
$ws.Range("N46").Value = -3961

$ws.Range("H100").Value = 3595.8235
$ws.Range("I100").Value = 3702.0667
$ws.Range("J100").Value = 2799
$ws.Range("K100").Value = 3702.0667
$ws.Range("L100").Value = 2799
$ws.Range("M100").Value = -3161.0667
$ws.Range("N100").Value = -3881

$ws.Range("H126").Value = 6617.3335
$ws.Range("I126").Value = 2877.8333
$ws.Range("K126").Value = 8633.499899999999
$ws.Range("M126").Value = -6163.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51144

$ws.Range("H122").Value = 7288.35
$ws.Range("I122").Value = 5979.077
$ws.Range("J122").Value = 9719.857
$ws.Range("K122").Value = 17937.231
$ws.Range("L122").Value = 29159.571
$ws.Range("M122").Value = -15487.231
$ws.Range("N122").Value = -34059.571

$ws.Range("H126").Value = 2788.3333
$ws.Range("J126").Value = 2699
$ws.Range("L126").Value = 8097
$ws.Range("N126").Value = -13037

$ws.Range("H136").Value = 736.46155
$ws.Range("I136").Value = 689
$ws.Range("K136").Value = 2067
$ws.Range("M136").Value = 483
